$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

# Row 8: copy the existing date-format style (from B2) onto B8 so the
# inserted date reuses the same cellXf instead of minting a new one.
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

# Row 8 & 9: copy the existing percentage-format style (from E2) onto
# E8/E9 for the same reason.
$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E9").PasteSpecial(-4122)

# Row 8 values
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "3/3/2022"
$ws.Cells.Item(8, 3).Value = "RPA GSS"
$ws.Cells.Item(8, 4).Value = "1. Invoice_generation has been completed, tested and it is running smoothly"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = "Completed"

# Row 9 values
$ws.Cells.Item(9, 4).Value = "2. Task of Service Order Pending  is work in progress"
$ws.Cells.Item(9, 5).Value = 0.1
$ws.Cells.Item(9, 6).Value = "WIP"

# Move the selection to D19, matching the author's final cursor position.
$null = $ws.Range("D19").Select()
